# "grades and budget update"
# Update the grade-tracking worksheet (Sheet1) with the latest scores:
#   - Quiz grade earned (E8) changed from 1 -> 0.95
#   - HW grade earned (D11) changed from 0.9 -> 1
# The percentage-available row (row 7) and the grade-total column (J)
# use formulas, so they recalculate automatically.
# Also move the active cell selection from D9 to E9, matching the
# cursor position left by the editor when the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E8").Value = 0.95
$ws.Range("D11").Value = 1

$ws.Activate()
$ws.Range("E9").Select()
